$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 1133, shifting existing rows 1133:1213 down to 1135:1215
$ws.Range("1133:1134").Insert()

# Fill new row 1133 with the new weekly record
$ws.Range("A1133").Value = 12
$ws.Range("B1133").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C1133").Value = "Metropolitana"
$ws.Range("D1133").Value = 44585
$ws.Range("E1133").Value = 13
$ws.Range("F1133").Value = 100112020
$ws.Range("G1133").Value = "Tomate"
$ws.Range("H1133").Value = "Larga vida"
$ws.Range("I1133").Value = "Extra"
$ws.Range("J1133").Value = 380
$ws.Range("K1133").Value = 13000
$ws.Range("L1133").Value = 13000
$ws.Range("M1133").Value = 13000
$ws.Range("N1133").Value = "$/bandeja 18 kilos"
$ws.Range("O1133").Value = "Provincia de Quillota"
$ws.Range("P1133").Value = 722
$ws.Range("Q1133").Value = 18
$ws.Range("R1133").Value = "Hortaliza"

# Fill new row 1134 with the new weekly record
$ws.Range("A1134").Value = 12
$ws.Range("B1134").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C1134").Value = "Metropolitana"
$ws.Range("D1134").Value = 44585
$ws.Range("E1134").Value = 13
$ws.Range("F1134").Value = 100112020
$ws.Range("G1134").Value = "Tomate"
$ws.Range("H1134").Value = "Larga vida"
$ws.Range("I1134").Value = "Primera"
$ws.Range("J1134").Value = 880
$ws.Range("K1134").Value = 9000
$ws.Range("L1134").Value = 11000
$ws.Range("M1134").Value = 9909
$ws.Range("N1134").Value = "$/bandeja 18 kilos"
$ws.Range("O1134").Value = "Provincia de Quillota"
$ws.Range("P1134").Value = 550
$ws.Range("Q1134").Value = 18
$ws.Range("R1134").Value = "Hortaliza"
